## ---------------------------------------------------------------------
## Adds the new "Hoja3" lookup sheet, appends two new asset rows
## (DINA T / AKACIAS) to Hoja1, tweaks a couple of view selections and
## one cell's number format, matching the authored diff.
## ---------------------------------------------------------------------

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws2 = $wb.Worksheets.Item("Hoja2")

## 1. F32 on Hoja1 gets an explicit "General" number format (keeps its border) -----
$ws1.Range("F32").NumberFormat = "General"

## 2. Append the two new rows (47 & 48) to Hoja1 -----------------------
# Row 47 - DINA T
$ws1.Range("A47").Value = "DINA T"
$ws1.Range("B47").Value = "activo"
$ws1.Range("E47").Value = "ECOPETROL"
$ws1.Range("F47").Value = 281
$ws1.Range("G47").Value = 74.25
$ws1.Range("H47").Formula = "=F47-G47"
$ws1.Range("I47").Formula = "=G47/F47"
$ws1.Range("J47").Value = 7.09
$ws1.Range("K47").Value = 0
$ws1.Range("L47").Value = 0
$ws1.Range("M47").Formula = "=(J47+K47+L47+N47+O47)/F47"
$ws1.Range("N47").Value = 11.78
$ws1.Range("O47").Value = 5.83
$ws1.Range("P47").Formula = "=K47+L47+N47+O47"
$ws1.Range("Q47").Value = 2809
$ws1.Range("R47").Value = 45082
$ws1.Range("S47").Value = 56.4
$ws1.Range("T47").Value = 4.5

# Row 48 - AKACIAS
$ws1.Range("A48").Value = "AKACIAS"
$ws1.Range("B48").Value = "activo"
$ws1.Range("E48").Value = "ECOPETROL"
$ws1.Range("F48").Value = 2471
$ws1.Range("G48").Value = 23
$ws1.Range("H48").Formula = "=F48-G48"
$ws1.Range("I48").Formula = "=G48/F48"
$ws1.Range("J48").Value = 74
$ws1.Range("K48").Value = 78
$ws1.Range("L48").Value = 52
$ws1.Range("M48").Formula = "=(J48+K48+L48+N48+O48)/F48"
$ws1.Range("N48").Value = 283
$ws1.Range("O48").Value = 0
$ws1.Range("P48").Formula = "=K48+L48+N48+O48"
$ws1.Range("Q48").Value = 16775.614657534246
$ws1.Range("R48").Value = 45263
$ws1.Range("S48").Value = 654
$ws1.Range("T48").Value = 4.5

## Apply date formatting to the two "date" cells --------------------
$ws1.Range("R47:R48").NumberFormat = "mm-dd-yy"

## Highlight the new rows in red, like the author did -----------------
$ws1.Range("A47:T48").Font.Color = 255
# F47 was left out of the highlight in the source workbook
$ws1.Range("F47").Font.Color = -4105

## 3. Add the "Hoja3" lookup sheet, placed after Hoja2 -----------------
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $last)
$ws3.Name = "Hoja3"

$ws3.Range("G2").Value = "ACTIVO"
$ws3.Range("H2").Value = "DINA T"
$ws3.Range("I2").Value = "FUENTE"

$ws3.Range("G3").Value = "GERENCIA VAS"

$ws3.Range("G4").Value = "Operador"
$ws3.Range("H4").Value = "ECP"

$ws3.Range("G5").Value = "OOIP - MBls"
$ws3.Range("H5").Value = 281
$ws3.Range("I5").Value = "OOIP Oficial"

$ws3.Range("G6").Value = "Np - MBls"
$ws3.Range("H6").Value = 74.25
$ws3.Range("I6").Value = "Dic/2021-OFM"

$ws3.Range("G7").Value = "FR_Act"
$ws3.Range("H7").Value = 26.4

$ws3.Range("G8").Value = "P1 - - MBls"
$ws3.Range("H8").Value = 7.09
$ws3.Range("I8").Value = "Balance de reservas oficial"

$ws3.Range("G9").Value = "P2 - MBls"
$ws3.Range("H9").Value = 0

$ws3.Range("G10").Value = "P3 - - MBls"
$ws3.Range("H10").Value = 0

$ws3.Range("G11").Value = "RC- MBls"
$ws3.Range("H11").Value = 11.78
$ws3.Range("I11").Value = "Balance de RC oficial"

$ws3.Range("G12").Value = "RV - MBls"
$ws3.Range("H12").Value = 5.83
$ws3.Range("I12").Value = "PLP 2021"

$ws3.Range("G13").Value = "Fecha_finalización contrato"
$ws3.Range("H13").Value = 45082
$ws3.Range("H13").NumberFormat = "mm-dd-yy"

$ws3.Range("G14").Value = "VPN activo (MUSD)"
$ws3.Range("H14").Value = 56.4
$ws3.Range("I14").Value = "Portafolio/2022"

## Column widths on Hoja3 ------------------------------------------------
$ws3.Columns.Item(7).ColumnWidth = 25.42578125
$ws3.Columns.Item(8).ColumnWidth = 15.140625
$ws3.Columns.Item(9).ColumnWidth = 24.5703125

## 4. View / selection tweaks -------------------------------------------
# Hoja2: selection moves to H23 (keep Hoja1 as the active tab afterwards)
$ws2.Activate()
$ws2.Range("H23").Select()

# Hoja1: freeze the header row, scroll down, select E49
$ws1.Activate()
$ws1.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws1.Range("E49").Select()

# Hoja3: leave the selection on H14
$ws3.Activate()
$ws3.Range("H14").Select()

$ws1.Activate()
